$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Baz " -> "Baz": strip the trailing space from the shared string used by A2
$ws.Range("A2").Value = "Baz"

# New thin border around B2 (creates border #1 + a cellXf that references it)
$ws.Range("B2").Borders.LineStyle = 1

# Active selection moves to B2
$ws.Range("B2").Select() | Out-Null
